# Auto-generated Excel COM-interop script to apply Phantom_Profits value updates
# Data source: cached market-board snapshot columns (H:N) across all 8 leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5245.4546
$ws.Range("I18").Value = 5455.5557
$ws.Range("K18").Value = 5455.5557
$ws.Range("M18").Value = -5171.5557
$ws.Range("H70").Value = 2274.5715
$ws.Range("J70").Value = 2799.1667
$ws.Range("L70").Value = 8397.500100000001
$ws.Range("N70").Value = -8937.500100000001
$ws.Range("H73").Value = 2274.5715
$ws.Range("J73").Value = 2799.1667
$ws.Range("L73").Value = 8397.500100000001
$ws.Range("N73").Value = -10269.5001
$ws.Range("H100").Value = 2956.5715
$ws.Range("I100").Value = 3382.8333
$ws.Range("J100").Value = 399
$ws.Range("K100").Value = 3382.8333
$ws.Range("L100").Value = 399
$ws.Range("M100").Value = -2841.8333
$ws.Range("N100").Value = -1481
$ws.Range("H138").Value = 2483.3618
$ws.Range("J138").Value = 2221.8823
$ws.Range("L138").Value = 6665.646900000001
$ws.Range("N138").Value = -16945.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 13555
$ws.Range("I56").Value = 7000
$ws.Range("J56").Value = 20110
$ws.Range("K56").Value = 7000
$ws.Range("L56").Value = 20110
$ws.Range("M56").Value = -6258
$ws.Range("N56").Value = -21594
$ws.Range("H74").Value = 2942.45
$ws.Range("I74").Value = 1402.1
$ws.Range("K74").Value = 1402.1
$ws.Range("M74").Value = -528.0999999999999
$ws.Range("H77").Value = 2942.45
$ws.Range("I77").Value = 1402.1
$ws.Range("K77").Value = 7010.5
$ws.Range("M77").Value = -2642.5
$ws.Range("H88").Value = 3048
$ws.Range("J88").Value = 3056
$ws.Range("L88").Value = 3056
$ws.Range("N88").Value = -3868
$ws.Range("H91").Value = 3048
$ws.Range("J91").Value = 3056
$ws.Range("L91").Value = 3056
$ws.Range("N91").Value = -5864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 21525
$ws.Range("I86").Value = 21850.5
$ws.Range("J86").Value = 20548.5
$ws.Range("K86").Value = 21850.5
$ws.Range("L86").Value = 20548.5
$ws.Range("M86").Value = -20727.5
$ws.Range("N86").Value = -22794.5
$ws.Range("H89").Value = 21525
$ws.Range("I89").Value = 21850.5
$ws.Range("J89").Value = 20548.5
$ws.Range("K89").Value = 109252.5
$ws.Range("L89").Value = 102742.5
$ws.Range("M89").Value = -103636.5
$ws.Range("N89").Value = -113974.5
$ws.Range("H94").Value = 1650.3572
$ws.Range("I94").Value = 1494.5
$ws.Range("K94").Value = 1494.5
$ws.Range("M94").Value = -1043.5
$ws.Range("H105").Value = 2044.5
$ws.Range("I105").Value = 2153.4
$ws.Range("K105").Value = 2153.4
$ws.Range("M105").Value = -406.4000000000001
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 1244.68
$ws.Range("I134").Value = 1213.2084
$ws.Range("K134").Value = 3639.6252
$ws.Range("M134").Value = -1104.6252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 994.75
$ws.Range("I16").Value = 994.75
$ws.Range("K16").Value = 994.75
$ws.Range("M16").Value = -707.75
$ws.Range("H51").Value = 39364.07
$ws.Range("I51").Value = 35999.875
$ws.Range("J51").Value = 43849.668
$ws.Range("K51").Value = 35999.875
$ws.Range("L51").Value = 43849.668
$ws.Range("M51").Value = -35263.875
$ws.Range("N51").Value = -45321.668
$ws.Range("H61").Value = 39364.07
$ws.Range("I61").Value = 35999.875
$ws.Range("J61").Value = 43849.668
$ws.Range("K61").Value = 35999.875
$ws.Range("L61").Value = 43849.668
$ws.Range("M61").Value = -35651.875
$ws.Range("N61").Value = -44545.668
$ws.Range("H62").Value = 4149.364
$ws.Range("J62").Value = 4327.3335
$ws.Range("L62").Value = 4327.3335
$ws.Range("N62").Value = -5575.3335
$ws.Range("H65").Value = 4149.364
$ws.Range("J65").Value = 4327.3335
$ws.Range("L65").Value = 21636.6675
$ws.Range("N65").Value = -27876.6675
$ws.Range("H74").Value = 50313.5
$ws.Range("J74").Value = 50313.5
$ws.Range("L74").Value = 50313.5
$ws.Range("N74").Value = -52061.5
$ws.Range("H77").Value = 50313.5
$ws.Range("J77").Value = 50313.5
$ws.Range("L77").Value = 150940.5
$ws.Range("N77").Value = -159676.5
$ws.Range("H107").Value = 1619.1428
$ws.Range("I107").Value = 1555.6666
$ws.Range("K107").Value = 1555.6666
$ws.Range("M107").Value = 364.3334
$ws.Range("H113").Value = 994.75
$ws.Range("I113").Value = 994.75
$ws.Range("K113").Value = 994.75
$ws.Range("M113").Value = 1175.25
$ws.Range("H122").Value = 2233.5
$ws.Range("J122").Value = 1998.3334
$ws.Range("L122").Value = 5995.0002
$ws.Range("N122").Value = -10895.0002
$ws.Range("H134").Value = 1633.75
$ws.Range("I134").Value = 1464.091
$ws.Range("K134").Value = 4392.272999999999
$ws.Range("M134").Value = -1857.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 62773.35
$ws.Range("J34").Value = 76147
$ws.Range("L34").Value = 228441
$ws.Range("N34").Value = -228609
$ws.Range("H39").Value = 7299.923
$ws.Range("J39").Value = 7283.25
$ws.Range("L39").Value = 21849.75
$ws.Range("N39").Value = -22437.75
$ws.Range("H46").Value = 1088.8
$ws.Range("J46").Value = 1111
$ws.Range("L46").Value = 3333
$ws.Range("N46").Value = -3515
$ws.Range("H55").Value = 300.8
$ws.Range("J55").Value = 399.5
$ws.Range("L55").Value = 1198.5
$ws.Range("N55").Value = -1552.5
$ws.Range("H116").Value = 2949
$ws.Range("J116").Value = 898
$ws.Range("L116").Value = 2694
$ws.Range("N116").Value = -9578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5013.7144
$ws.Range("J80").Value = 6019.4
$ws.Range("L80").Value = 6019.4
$ws.Range("N80").Value = -8015.4
$ws.Range("H83").Value = 5013.7144
$ws.Range("J83").Value = 6019.4
$ws.Range("L83").Value = 30097
$ws.Range("N83").Value = -40081
$ws.Range("H102").Value = 1931.4762
$ws.Range("I102").Value = 1592.6316
$ws.Range("K102").Value = 1592.6316
$ws.Range("M102").Value = 29.36840000000007
$ws.Range("H107").Value = 3287.111
$ws.Range("I107").Value = 617.4
$ws.Range("J107").Value = 6624.25
$ws.Range("K107").Value = 617.4
$ws.Range("L107").Value = 6624.25
$ws.Range("M107").Value = 1302.6
$ws.Range("N107").Value = -10464.25
$ws.Range("H133").Value = 78999.39999999999
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 78999.39999999999
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 78999.39999999999
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -89119.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 812.3333
$ws.Range("I55").Value = 400.15384
$ws.Range("J55").Value = 1299.4546
$ws.Range("K55").Value = 400.15384
$ws.Range("L55").Value = 1299.4546
$ws.Range("M55").Value = -227.15384
$ws.Range("N55").Value = -1645.4546
$ws.Range("H61").Value = 3758
$ws.Range("I61").Value = 3758
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3758
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3556
$ws.Range("N61").ClearContents()
$ws.Range("H93").Value = 926.5
$ws.Range("I93").Value = 847.25
$ws.Range("J93").Value = 1085
$ws.Range("K93").Value = 847.25
$ws.Range("L93").Value = 1085
$ws.Range("M93").Value = 400.75
$ws.Range("N93").Value = -3581
$ws.Range("H108").Value = 84682.664
$ws.Range("J108").Value = 84682.664
$ws.Range("L108").Value = 84682.664
$ws.Range("N108").Value = -92362.664
$ws.Range("H113").Value = 3758
$ws.Range("I113").Value = 3758
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3758
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1588
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 201016.33
$ws.Range("I61").Value = 201016.33
$ws.Range("K61").Value = 201016.33
$ws.Range("M61").Value = -200724.33
$ws.Range("H125").Value = 60499
$ws.Range("J125").Value = 60499
$ws.Range("L125").Value = 60499
$ws.Range("N125").Value = -70339
$ws.Range("H126").Value = 899.6667
$ws.Range("J126").Value = 598.5
$ws.Range("L126").Value = 1795.5
$ws.Range("N126").Value = -6735.5
$ws.Range("H132").Value = 333338000
$ws.Range("I132").Value = 6999
$ws.Range("K132").Value = 20997
$ws.Range("M132").Value = -18467
$ws.Range("H138").Value = 73665.336
$ws.Range("J138").Value = 73665.336
$ws.Range("L138").Value = 73665.336
$ws.Range("N138").Value = -83945.336
